$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "incorrect"
$ws.Range("C2").Value = 0
$ws.Range("F2").Value = 1

$ws.Range("C3").Value = 0.6665308908294946
$ws.Range("F3").Value = 0.3334691091705055

$ws.Range("B5").Value = "incorrect"
$ws.Range("C5").Value = 0.3151691432172188
$ws.Range("D5").Value = 0.6848308567827812

$ws.Range("B6").Value = "incorrect"
$ws.Range("C6").Value = 0.3264818596618402
$ws.Range("D6").Value = 0.3273263381045693
$ws.Range("E6").Value = 0.3461918022335906

$ws.Range("B7").Value = "incorrect"
$ws.Range("C7").Value = 0.3435351089752017
$ws.Range("F7").Value = 0.6564648910247981

$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 0
$ws.Range("F8").Value = 0

$ws.Range("B10").Value = "incorrect"
$ws.Range("C10").Value = 0.6669168179794867
$ws.Range("D10").Value = 0
$ws.Range("F10").Value = 0.3330831820205133

$ws.Range("B11").Value = "correct"
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 1

$ws.Range("B12").Value = "correct"
$ws.Range("C12").Value = 0.3264689577657587
$ws.Range("D12").Value = 0.3421267112218081
$ws.Range("F12").Value = 0.3314043310124331

$ws.Range("D13").Value = 0.6704549760390368
$ws.Range("F13").Value = 0.3295450239609632

$ws.Range("B14").Value = "incorrect"
$ws.Range("C14").Value = 0.6782308396048843
$ws.Range("D14").Value = 0.3217691603951158
$ws.Range("E14").Value = 0

$ws.Range("B15").Value = "incorrect"
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 0

$ws.Range("B16").Value = "incorrect"
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0

$ws.Range("B17").Value = "incorrect"
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 0

$ws.Range("B18").Value = "incorrect"
$ws.Range("C18").Value = 1
$ws.Range("E18").Value = 0

$ws.Range("B19").Value = "incorrect"
$ws.Range("D19").Value = 0.6745626536237902
$ws.Range("E19").Value = 0.3254373463762097

$ws.Range("C20").Value = 0.3511738586968672
$ws.Range("F20").Value = 0.6488261413031329

$ws.Range("E21").Value = 0.334688955108254
$ws.Range("F21").Value = 0.6653110448917461

$ws.Range("C22").Value = 0.6698870206226143
$ws.Range("D22").Value = 0.3301129793773858
$ws.Range("E22").Value = 0

$ws.Range("B23").Value = "incorrect"
$ws.Range("C23").Value = 0.3578590534864374
$ws.Range("D23").Value = 0.6421409465135626
$ws.Range("F23").Value = 0

$ws.Range("C24").Value = 1
$ws.Range("D24").Value = 0

$ws.Range("C25").Value = 0
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 0

$ws.Range("C26").Value = 0.3275421410722453
$ws.Range("E26").Value = 0.3365458514056789
$ws.Range("F26").Value = 0.3359120075220757
